$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    # Force the value to be written/stay as literal text (not auto-converted
    # to a number/date by Excel), matching the source inline-string cells,
    # then strip the temporary number-format override so no stray style sticks.
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).ClearFormats()
}

Set-TextValue $ws "D2" "69.339.20"
$ws.Range("E2").Value = "  +2.45%  "
Set-TextValue $ws "D3" "3.386.14"
$ws.Range("E3").Value = "  +1.53%  "
Set-TextValue $ws "D4" "0.999"
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws "D5" "588.70"
$ws.Range("E5").Value = "  +1.38%  "
Set-TextValue $ws "D6" "180.84"
$ws.Range("E6").Value = "  +2.89%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("E9").Value = "  +8.94%  "
$ws.Range("E10").Value = "  +1.43%  "
Set-TextValue $ws "D11" "48.97"
$ws.Range("E11").Value = "  +5.82%  "
Set-TextValue $ws "D12" "0.0000284"
$ws.Range("E12").Value = "  +4.63%  "
Set-TextValue $ws "D13" "682.59"
$ws.Range("E13").Value = "  -2.86%  "
Set-TextValue $ws "D14" "8.63"
$ws.Range("E14").Value = "  +2.23%  "
Set-TextValue $ws "D15" "3.930.58"
$ws.Range("E15").Value = "  +1.22%  "
Set-TextValue $ws "D16" "69.378.47"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("E17").Value = "  +1.69%  "
Set-TextValue $ws "D18" "3.375.47"
$ws.Range("E18").Value = "  +1.04%  "
Set-TextValue $ws "D19" "17.79"
$ws.Range("E19").Value = "  +2.35%  "
Set-TextValue $ws "D20" "11.42"
$ws.Range("E20").Value = "  +3.90%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  -0.32%  "
Set-TextValue $ws "D23" "17.11"
Set-TextValue $ws "D24" "104.79"
$ws.Range("E24").Value = "  +6.58%  "
Set-TextValue $ws "D25" "3.95"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  +1.52%  "
Set-TextValue $ws "D27" "9.65"
$ws.Range("E27").Value = "  +2.14%  "
Set-TextValue $ws "D28" "34.39"
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("E29").Value = "  +2.10%  "
Set-TextValue $ws "D30" "7.03"
$ws.Range("E30").Value = "  -1.23%  "
Set-TextValue $ws "D31" "11.21"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D32" "558.87"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D33" "3.64"
$ws.Range("E33").Value = "  +9.41%  "
Set-TextValue $ws "D34" "0.107"
$ws.Range("E34").Value = "  +1.32%  "
Set-TextValue $ws "D35" "58.11"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D36" "1.00"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws "D37" "3.725.57"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("E38").Value = "  +8.36%  "
Set-TextValue $ws "D39" "35.06"
$ws.Range("E39").Value = "  +3.00%  "
$ws.Range("E40").Value = "  +2.07%  "
Set-TextValue $ws "D41" "0.0₃0708"
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("E43").Value = "  +1.34%  "
Set-TextValue $ws "D44" "0.0419"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("E47").Value = "  +1.33%  "
Set-TextValue $ws "D48" "1.40"
$ws.Range("E48").Value = "  +5.72%  "
$ws.Range("E49").Value = "  -0.11%  "
Set-TextValue $ws "D50" "132.69"
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("E51").Value = "  -3.44%  "
